$wb = $excel.ActiveWorkbook

# --- Sheet2 ("Sheet2 - Numbers"): add a new column AA with values 100..129 ---
$ws2 = $wb.Worksheets.Item(2)

for ($r = 1; $r -le 30; $r++) {
    $ws2.Cells.Item($r, 27).Value = 99 + $r   # column 27 = AA
}

# Make Sheet2 the active sheet/tab and select AA1:AA30 (matches the new
# sheetView selection: activeCell AA1, sqref AA1:AA30, tabSelected).
$ws2.Select()
$ws2.Range("AA1:AA30").Select()

# --- Sheet4 ("Sheet4 - Dates"): page setup paper size changed from 0 to 9 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.PageSetup.PaperSize = 9
